$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Establish new shared strings in the order they appear in the target file:
# index 8 = "Passwort verschlüsseln (wenigstens MD5 Hash)"
# index 9 = "erledigr"
# index 10 = "erledigt"
$ws.Range("A7").Value = "Passwort verschlüsseln (wenigstens MD5 Hash)"
$ws.Range("E7").Value = "erledigr"
$ws.Range("E3").Value = "erledigt"

# Row 4: add E4 = "erledigt" (reuses existing shared string)
$ws.Range("E4").Value = "erledigt"

$green = 5287936   # RGB(0,176,80) -> OLE BGR value

# Apply green font to the affected rows
$ws.Range("A3:E3").Font.Color = $green
$ws.Range("A4:E4").Font.Color = $green
$ws.Range("A7:E7").Font.Color = $green

# Update selection to A2
$ws.Range("A2").Select()

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
